$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column B (rows 2-67)
$bValues = @(
    -0.2998608192475538,
    0.2222994047451824,
    0.1596690214676689,
    -0.05157580157622221,
    0.1524653381963466,
    -0.4397818082765901,
    -0.2891752268856903,
    -0.4046922839984844,
    0.3269025025094837,
    -0.2058413461554524,
    -0.08305774230531343,
    -0.001980202257975946,
    0.03918117179674931,
    -0.07737564180358925,
    0.2727842049609036,
    0.5925036014333478,
    0.0640445261876246,
    0.3982653262601039,
    0.2487168247563672,
    0.4435913990058672,
    0.3552801054998972,
    -0.06598158323609662,
    4.64652098450761,
    0.4943215704727376,
    0.3944919096527825,
    0.3478665255842199,
    1.131311382761552,
    5.989697285148654,
    1.040980659047581,
    -0.146546836024889,
    0.7828745218074415,
    0.9086159910027884,
    -0.6322938977899302,
    0.8308202261936922,
    0.7845997278769563,
    0.7673247537313208,
    0.7614248359261103,
    0.5801498560662414,
    0.7500862168053677,
    0.5740292953606498,
    0.7219757400546358,
    0.7344036898121714,
    0.6724262516180799,
    0.6557132659512696,
    -1.259310717349668,
    -0.9751540048435037,
    -0.8699514465962723,
    -0.6376377267191787,
    -0.05228625112963124,
    -0.8577452128631521,
    -0.8577452128631521,
    -1.085572194885408,
    -0.191313546284351,
    -0.9973543927928555,
    -0.8989397862939991,
    -0.9697356462549301,
    -1.173013568950063,
    -0.8758385668749407,
    -0.5276024229975791,
    0.3647378274819701,
    -1.266473933546654,
    -0.7710400446527043,
    -0.9064267757069395,
    -0.1456098824911458,
    -0.8427701866883508,
    -0.8447519046281721
)

# New values for column C (rows 2-67)
$cValues = @(
    -0.2844306209533958,
    0.2323567825548751,
    0.1730073447847666,
    -0.04133725088963019,
    0.1693687572142441,
    -0.4288345250936945,
    -0.2655955976075189,
    -0.3945937732805802,
    0.3446233488301431,
    -0.1976879972501496,
    -0.05729188132265255,
    0.0006439044102356604,
    0.079516501932351,
    -0.02761333128479813,
    0.3387468086258007,
    0.6211163730947561,
    0.04370221162946492,
    0.4078726304003948,
    0.3079483271226175,
    0.5110900338356539,
    0.3917780878624746,
    -0.03362055528175069,
    4.657238259679452,
    0.4815387829650647,
    0.36244794794862,
    0.3125625161631608,
    1.073841271021365,
    5.507915010888566,
    0.976185683296799,
    -0.2171781629357534,
    0.7576030272578753,
    0.893840682684764,
    -0.6424596611186806,
    0.8211784480373452,
    0.7682661920851495,
    0.7468983712376456,
    0.7381370088434778,
    0.5785019198624513,
    0.7503044209841752,
    0.5657034547119154,
    0.6917388954469098,
    0.7183746850466973,
    0.6655174257823021,
    0.6393301613064649,
    -1.262619605060654,
    -0.9794514840791761,
    -0.8748889017938992,
    -0.6400733031808122,
    -0.05225096410204793,
    -0.8604280593777454,
    -0.8604280593777454,
    -1.09905820809552,
    -0.1888775749458106,
    -0.9998245551834358,
    -0.8919475490897619,
    -0.9519109274952845,
    -1.140645740446841,
    -0.857186425265752,
    -0.5000282140209131,
    0.3684809577194981,
    -1.240715157239431,
    -0.7331895034939727,
    -0.8967893644422141,
    -0.1189287952192628,
    -0.8091611560749737,
    -0.8005741458551142
)

$startRow = 2
for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
}